$d = $word.ActiveDocument
$p1 = $d.Paragraphs.Item(1)
$r = $p1.Range
$frag = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:v="urn:schemas-microsoft-com:vml" xmlns:o="urn:schemas-microsoft-com:office:office" xmlns:w10="urn:schemas-microsoft-com:office:word" w:rsidR="005A0901" w:rsidRDefault="005A0901" w:rsidP="005A0901"><w:r><w:rPr><w:noProof/><w:lang w:eastAsia="es-ES"/></w:rPr><w:pict><v:group id="_x0000_s1041" style="position:absolute;margin-left:-5.2pt;margin-top:-16.05pt;width:611.8pt;height:649.85pt;z-index:251662336;mso-position-horizontal-relative:page;mso-position-vertical-relative:margin" coordorigin=",1440" coordsize="12239,12960" o:allowincell="f"><v:group id="_x0000_s1042" style="position:absolute;top:9661;width:12239;height:4739;mso-position-horizontal:center;mso-position-horizontal-relative:margin;mso-position-vertical:bottom;mso-position-vertical-relative:margin" coordorigin="-6,3399" coordsize="12197,4253"><v:group id="_x0000_s1043" style="position:absolute;left:-6;top:3717;width:12189;height:3550" coordorigin="18,7468" coordsize="12189,3550"><v:shape id="_x0000_s1044" style="position:absolute;left:18;top:7837;width:7132;height:2863" coordsize="7132,2863" path="m,l17,2863,7132,2578r,-2378l,xe" fillcolor="#a7bfde" stroked="f"><v:fill opacity=".5"/><v:path arrowok="t"/></v:shape><v:shape id="_x0000_s1045" style="position:absolute;left:7150;top:7468;width:3466;height:3550" coordsize="3466,3550" path="m,569l,2930r3466,620l3466,,,569xe" fillcolor="#d3dfee" stroked="f"><v:fill opacity=".5"/><v:path arrowok="t"/></v:shape><v:shape id="_x0000_s1046" style="position:absolute;left:10616;top:7468;width:1591;height:3550" coordsize="1591,3550" path="m,l,3550,1591,2746r,-2009l,xe" fillcolor="#a7bfde" stroked="f"><v:fill opacity=".5"/><v:path arrowok="t"/></v:shape></v:group><v:shape id="_x0000_s1047" style="position:absolute;left:8071;top:4069;width:4120;height:2913" coordsize="4120,2913" path="m1,251l,2662r4120,251l4120,,1,251xe" fillcolor="#d8d8d8" stroked="f"><v:path arrowok="t"/></v:shape><v:shape id="_x0000_s1048" style="position:absolute;left:4104;top:3399;width:3985;height:4236" coordsize="3985,4236" path="m,l,4236,3985,3349r,-2428l,xe" fillcolor="#bfbfbf" stroked="f"><v:path arrowok="t"/></v:shape><v:shape id="_x0000_s1049" style="position:absolute;left:18;top:3399;width:4086;height:4253" coordsize="4086,4253" path="m4086,r-2,4253l,3198,,1072,4086,xe" fillcolor="#d8d8d8" stroked="f"><v:path arrowok="t"/></v:shape><v:shape id="_x0000_s1050" style="position:absolute;left:17;top:3617;width:2076;height:3851" coordsize="2076,3851" path="m,921l2060,r16,3851l,2981,,921xe" fillcolor="#d3dfee" stroked="f"><v:fill opacity="45875f"/><v:path arrowok="t"/></v:shape><v:shape id="_x0000_s1051" style="position:absolute;left:2077;top:3617;width:6011;height:3835" coordsize="6011,3835" path="m,l17,3835,6011,2629r,-1390l,xe" fillcolor="#a7bfde" stroked="f"><v:fill opacity="45875f"/><v:path arrowok="t"/></v:shape><v:shape id="_x0000_s1052" style="position:absolute;left:8088;top:3835;width:4102;height:3432" coordsize="4102,3432" path="m,1038l,2411,4102,3432,4102,,,1038xe" fillcolor="#d3dfee" stroked="f"><v:fill opacity="45875f"/><v:path arrowok="t"/></v:shape></v:group><v:rect id="_x0000_s1053" style="position:absolute;left:1800;top:1440;width:8638;height:1935;mso-position-horizontal:center;mso-position-horizontal-relative:margin;mso-position-vertical:top;mso-position-vertical-relative:margin" filled="f" stroked="f"><v:textbox style="mso-next-textbox:#_x0000_s1053;mso-fit-shape-to-text:t"><w:txbxContent><w:p w:rsidR="005A0901" w:rsidRPr="00D92AA9" w:rsidRDefault="005A0901" w:rsidP="005A0901"><w:pPr><w:spacing w:after="0"/><w:jc w:val="center"/><w:rPr><w:b/><w:bCs/><w:color w:val="808080"/><w:sz w:val="48"/><w:szCs w:val="48"/></w:rPr></w:pPr><w:r w:rsidRPr="00D92AA9"><w:rPr><w:b/><w:bCs/><w:color w:val="808080"/><w:sz w:val="48"/><w:szCs w:val="48"/></w:rPr><w:t>UNIVERSIDAD TECONOLÓGICA NACIONAL                    FACULTAD REGIONAL CÓRDOBA</w:t></w:r></w:p><w:p w:rsidR="005A0901" w:rsidRDefault="005A0901" w:rsidP="005A0901"><w:pPr><w:spacing w:after="0"/><w:rPr><w:b/><w:bCs/><w:color w:val="808080"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr></w:pPr></w:p></w:txbxContent></v:textbox></v:rect><v:rect id="_x0000_s1054" style="position:absolute;left:6494;top:11161;width:4998;height:1127;mso-position-horizontal-relative:margin;mso-position-vertical-relative:margin" filled="f" stroked="f"><v:textbox style="mso-next-textbox:#_x0000_s1054;mso-fit-shape-to-text:t"><w:txbxContent><w:p w:rsidR="005A0901" w:rsidRDefault="005A0901" w:rsidP="005A0901"><w:pPr><w:jc w:val="center"/><w:rPr><w:sz w:val="56"/><w:szCs w:val="56"/></w:rPr></w:pPr></w:p></w:txbxContent></v:textbox></v:rect><v:rect id="_x0000_s1055" style="position:absolute;left:1800;top:2294;width:8638;height:7268;mso-position-horizontal:center;mso-position-horizontal-relative:margin;mso-position-vertical-relative:margin;v-text-anchor:bottom" filled="f" stroked="f"><v:textbox style="mso-next-textbox:#_x0000_s1055"><w:txbxContent><w:p w:rsidR="005A0901" w:rsidRDefault="005A0901" w:rsidP="005A0901"><w:pPr><w:spacing w:after="0"/><w:jc w:val="center"/><w:rPr><w:b/><w:bCs/><w:color w:val="1F497D"/><w:sz w:val="72"/><w:szCs w:val="72"/></w:rPr></w:pPr></w:p><w:p w:rsidR="005A0901" w:rsidRDefault="005A0901" w:rsidP="005A0901"><w:pPr><w:spacing w:after="0"/><w:jc w:val="center"/><w:rPr><w:b/><w:bCs/><w:color w:val="1F497D"/><w:sz w:val="72"/><w:szCs w:val="72"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:color w:val="1F497D"/><w:sz w:val="72"/><w:szCs w:val="72"/></w:rPr><w:t>Proyecto Final               Optical Marketing</w:t></w:r></w:p><w:p w:rsidR="005A0901" w:rsidRDefault="005A0901" w:rsidP="005A0901"><w:pPr><w:jc w:val="center"/><w:rPr><w:b/><w:bCs/><w:color w:val="808080"/><w:sz w:val="44"/><w:szCs w:val="44"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:color w:val="808080"/><w:sz w:val="44"/><w:szCs w:val="44"/></w:rPr><w:t xml:space="preserve">Manual de Instalación </w:t></w:r></w:p><w:p w:rsidR="005A0901" w:rsidRPr="00D92AA9" w:rsidRDefault="005A0901" w:rsidP="005A0901"><w:pPr><w:spacing w:after="0"/><w:rPr><w:b/><w:bCs/><w:color w:val="808080"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:u w:val="single"/></w:rPr></w:pPr><w:r w:rsidRPr="00D92AA9"><w:rPr><w:b/><w:bCs/><w:color w:val="808080"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:u w:val="single"/></w:rPr><w:t>Profesores:</w:t></w:r></w:p><w:p w:rsidR="005A0901" w:rsidRDefault="005A0901" w:rsidP="005A0901"><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:b/><w:bCs/><w:color w:val="808080"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:r w:rsidRPr="00D92AA9"><w:rPr><w:b/><w:bCs/><w:color w:val="808080"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve">Zohil, Julio Cesar Nelson </w:t></w:r></w:p><w:p w:rsidR="005A0901" w:rsidRDefault="005A0901" w:rsidP="005A0901"><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:b/><w:bCs/><w:color w:val="808080"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:r w:rsidRPr="00D92AA9"><w:rPr><w:b/><w:bCs/><w:color w:val="808080"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve">Aquino, Francisco </w:t></w:r></w:p><w:p w:rsidR="005A0901" w:rsidRDefault="005A0901" w:rsidP="005A0901"><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:b/><w:bCs/><w:color w:val="808080"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:r w:rsidRPr="00D92AA9"><w:rPr><w:b/><w:bCs/><w:color w:val="808080"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>Jaime, Natalia</w:t></w:r></w:p><w:p w:rsidR="005A0901" w:rsidRPr="00D92AA9" w:rsidRDefault="005A0901" w:rsidP="005A0901"><w:pPr><w:spacing w:after="0"/><w:rPr><w:b/><w:bCs/><w:color w:val="808080"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:u w:val="single"/></w:rPr></w:pPr><w:r w:rsidRPr="00D92AA9"><w:rPr><w:b/><w:bCs/><w:color w:val="808080"/><w:sz w:val="28"/><w:szCs w:val="28"/><w:u w:val="single"/></w:rPr><w:t xml:space="preserve">Grupo 4:                                                                             </w:t></w:r></w:p><w:p w:rsidR="005A0901" w:rsidRPr="00D92AA9" w:rsidRDefault="005A0901" w:rsidP="005A0901"><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:b/><w:bCs/><w:color w:val="808080"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:r w:rsidRPr="00D92AA9"><w:rPr><w:b/><w:bCs/><w:color w:val="808080"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>Carlos Kapica 51482                                                                                       Rodrigo Liberal 51658</w:t></w:r></w:p><w:p w:rsidR="005A0901" w:rsidRPr="00D92AA9" w:rsidRDefault="005A0901" w:rsidP="005A0901"><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:b/><w:bCs/><w:color w:val="808080"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:r w:rsidRPr="00D92AA9"><w:rPr><w:b/><w:bCs/><w:color w:val="808080"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>Julián Peker 51395</w:t></w:r></w:p><w:p w:rsidR="005A0901" w:rsidRPr="00D92AA9" w:rsidRDefault="005A0901" w:rsidP="005A0901"><w:pPr><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:b/><w:bCs/><w:color w:val="808080"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr></w:pPr><w:r w:rsidRPr="00D92AA9"><w:rPr><w:b/><w:bCs/><w:color w:val="808080" w:themeColor="text1" w:themeTint="7F"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>Fernández David 53063</w:t></w:r></w:p></w:txbxContent></v:textbox></v:rect><w10:wrap anchorx="page" anchory="margin"/></v:group></w:pict></w:r></w:p>'
$r.InsertXML($frag)
Write-Output "done"
